# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$wsRush = $wb.Worksheets.Item("Rushing")

$wsRush.Range("E2").Value = 2

$wsRush.Range("C3").Value = 35
$wsRush.Range("D3").Value = 18
$wsRush.Range("E3").Value = 8

$wsRush.Range("C4").Value = 57
$wsRush.Range("D4").Value = 44
$wsRush.Range("E4").Value = 12

# --- Receiving sheet ---
$wsRecv = $wb.Worksheets.Item("Receiving")

$wsRecv.Range("C2").Value = 31
$wsRecv.Range("D2").Value = 23

$wsRecv.Range("C3").Value = 34
$wsRecv.Range("D3").Value = 31

$wsRecv.Range("C6").Value = 42
$wsRecv.Range("D6").Value = 30
$wsRecv.Range("E6").Value = 15
$wsRecv.Range("F6").Value = 7
$wsRecv.Range("G6").Value = 8

$wsRecv.Range("C7").Value = 55
$wsRecv.Range("D7").Value = 34
$wsRecv.Range("E7").Value = 11
$wsRecv.Range("G7").Value = 12

$wsRecv.Range("C8").Value = 41
$wsRecv.Range("D8").Value = 27
$wsRecv.Range("E8").Value = 16

$wsRecv.Range("C10").Value = 20
$wsRecv.Range("D10").Value = 11

$wsRecv.Range("C13").Value = 8
$wsRecv.Range("G13").Value = 2

$wsRecv.Range("C14").Value = 57
$wsRecv.Range("D14").Value = 38

$wsRecv.Range("C15").Value = 27
$wsRecv.Range("D15").Value = 19
